$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the name field (B2) and email field (D2) while keeping other
# row-2 values (C2, E2, F2) the same.
$ws.Range("B2").Value = "Priyaspr1"
$ws.Range("D2").Value = "SPRPriyas@gmail.com"

# Move the active selection from F2 to D2
$ws.Range("D2").Select()
